# Word COM-interop script implementing the "raven.docx" commit:
#
#   1. The first paragraph ("This is a Microsoft word document.") gets two
#      trailing spaces appended to its existing run, followed by a new
#      red-colored (FF0000) annotation reading
#        "(This is a change – Version for main branch)"
#      typed/applied as three separate runs.
#
#   2. The very last paragraph of the document -- the stray
#      "ank God almighty, we are free at last." line that trails the
#      poem's final "Shall be lifted—nevermore!" paragraph -- is deleted
#      entirely (including its paragraph mark), so the poem's last line
#      again sits directly above the section break.

$d = $word.ActiveDocument

# --- 1. Edit the first paragraph -----------------------------------------

$firstPara = $d.Paragraphs.Item(1).Range

# Position immediately before paragraph 1's end-of-paragraph mark.
$insertAt = $firstPara.End - 1
$d.Range($insertAt, $insertAt).InsertAfter("  ")

# Run 1: "(This is a change – Ve"  (en dash = U+2013)
$para = $d.Paragraphs.Item(1).Range
$run1Start = $para.End - 1
$d.Range($run1Start, $run1Start).InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$para = $d.Paragraphs.Item(1).Range
$run1End = $para.End - 1
$d.Range($run1Start, $run1End).Font.Color = 255

# Run 2: "rsion for main branch"
$run2Start = $run1End
$d.Range($run2Start, $run2Start).InsertAfter("rsion for main branch")
$para = $d.Paragraphs.Item(1).Range
$run2End = $para.End - 1
$d.Range($run2Start, $run2End).Font.Color = 255

# Run 3: ")"
$run3Start = $run2End
$d.Range($run3Start, $run3Start).InsertAfter(")")
$para = $d.Paragraphs.Item(1).Range
$run3End = $para.End - 1
$d.Range($run3Start, $run3End).Font.Color = 255

# --- 2. Remove the trailing "ank God almighty..." paragraph --------------

$target = $d.Content
$found = $target.Find.Execute("ank God almighty, we are free at last.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    [void]$target.Expand(4)  # wdParagraph: grow the hit to its whole paragraph (incl. mark)
    $target.Delete()
} else {
    # Fallback: the stray paragraph is always the document's last paragraph.
    $d.Paragraphs.Item($d.Paragraphs.Count).Range.Delete()
}
